# collect data thuong phat
# The "last_edited_time" column (D) for the rows that were re-synced on
# 2024-08-03 is consolidated onto a single, later timestamp
# (2024-08-03T03:54:00.000Z) instead of the two slightly different
# timestamps that were there before (03:28 and 03:29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2024-08-03T03:54:00.000Z"

$rows = @(4, 5, 6, 7, 8, 12, 13)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = $newTimestamp
}
